$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextCell {
    param($ws, $addr, $val)
    $cell = $ws.Range($addr)
    $cell.NumberFormat = "@"
    $cell.Value = $val
    $cell.NumberFormat = "General"
    $cell.Style = "Normal"
}

Set-TextCell $ws "D2" "70.699.04"
Set-TextCell $ws "E2" "  +1.95%  "
Set-TextCell $ws "D3" "3.473.44"
Set-TextCell $ws "E3" "  +2.55%  "
Set-TextCell $ws "D4" "1.00"
Set-TextCell $ws "E4" "  +0.11%  "
Set-TextCell $ws "D5" "588.38"
Set-TextCell $ws "E5" "  +0.28%  "
Set-TextCell $ws "D6" "180.73"
Set-TextCell $ws "E6" "  +0.56%  "
Set-TextCell $ws "D7" "3.463.61"
Set-TextCell $ws "E7" "  +2.40%  "
Set-TextCell $ws "D8" "0.602"
Set-TextCell $ws "E8" "  +1.11%  "
Set-TextCell $ws "E9" "  -0.07%  "
Set-TextCell $ws "D10" "0.206"
Set-TextCell $ws "E10" "  +5.93%  "
Set-TextCell $ws "D11" "0.596"
Set-TextCell $ws "E11" "  +0.92%  "
Set-TextCell $ws "D12" "49.56"
Set-TextCell $ws "E12" "  +2.27%  "
Set-TextCell $ws "D13" "0.0000287"
Set-TextCell $ws "E13" "  +2.13%  "
Set-TextCell $ws "D14" "695.18"
Set-TextCell $ws "E14" "  +2.76%  "
Set-TextCell $ws "D15" "8.80"
Set-TextCell $ws "E15" "  +2.21%  "
Set-TextCell $ws "D16" "4.027.86"
Set-TextCell $ws "E16" "  +2.38%  "
Set-TextCell $ws "D17" "70.767.71"
Set-TextCell $ws "E17" "  +1.94%  "
Set-TextCell $ws "D18" "3.474.58"
Set-TextCell $ws "E18" "  +1.08%  "
Set-TextCell $ws "E19" "  +1.11%  "
Set-TextCell $ws "D20" "17.96"
Set-TextCell $ws "E20" "  +2.00%  "
Set-TextCell $ws "D21" "11.51"
Set-TextCell $ws "E21" "  +2.16%  "
Set-TextCell $ws "D22" "0.917"
Set-TextCell $ws "E22" "  +1.62%  "
Set-TextCell $ws "D23" "5.48"
Set-TextCell $ws "E23" "  +0.89%  "
Set-TextCell $ws "D24" "17.23"
Set-TextCell $ws "E24" "  +0.50%  "
Set-TextCell $ws "D25" "102.05"
Set-TextCell $ws "E25" "  -1.38%  "
Set-TextCell $ws "E26" "  +1.16%  "
Set-TextCell $ws "E27" "  -0.12%  "
Set-TextCell $ws "D28" "9.75"
Set-TextCell $ws "E28" "  +0.96%  "
Set-TextCell $ws "D29" "34.04"
Set-TextCell $ws "E29" "  -0.06%  "
Set-TextCell $ws "D30" "8.92"
Set-TextCell $ws "E30" "  +2.53%  "
Set-TextCell $ws "E31" "  +3.91%  "
Set-TextCell $ws "D32" "3.96"
Set-TextCell $ws "E32" "  +10.32%  "
Set-TextCell $ws "D33" "579.25"
Set-TextCell $ws "E33" "  +4.25%  "
Set-TextCell $ws "E34" "  +0.49%  "
Set-TextCell $ws "D35" "59.11"
Set-TextCell $ws "E35" "  +1.74%  "
Set-TextCell $ws "E36" "  -1.88%  "
Set-TextCell $ws "D37" "1.00"
Set-TextCell $ws "E37" "  +0.03%  "
Set-TextCell $ws "D38" "3.608.94"
Set-TextCell $ws "E39" "  +1.57%  "
Set-TextCell $ws "D40" "35.74"
Set-TextCell $ws "E40" "  +1.71%  "
Set-TextCell $ws "B41" "Stacks"
Set-TextCell $ws "C41" "https://coinranking.com/coin/mMPrMcB7+stacks-stx"
Set-TextCell $ws "D41" "3.41"
Set-TextCell $ws "E41" "  +4.16%  "
Set-TextCell $ws "B42" "PEPE"
Set-TextCell $ws "C42" "https://coinranking.com/coin/03WI8NQPF+pepe-pepe"
Set-TextCell $ws "D42" "0.0₃0741"
Set-TextCell $ws "E42" "  +6.19%  "
Set-TextCell $ws "D43" "2.75"
Set-TextCell $ws "E43" "  +2.43%  "
Set-TextCell $ws "D44" "0.342"
Set-TextCell $ws "E44" "  +1.01%  "
Set-TextCell $ws "B45" "VeChain"
Set-TextCell $ws "C45" "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
Set-TextCell $ws "D45" "0.0430"
Set-TextCell $ws "E45" "  +1.89%  "
Set-TextCell $ws "B46" "ApeXProtocol"
Set-TextCell $ws "C46" "https://coinranking.com/coin/ze0N2Rcyu+apexprotocol-apex"
Set-TextCell $ws "D46" "3.37"
Set-TextCell $ws "E46" "  +2.95%  "
Set-TextCell $ws "D47" "2.73"
Set-TextCell $ws "E47" "  +1.95%  "
Set-TextCell $ws "D48" "1.46"
Set-TextCell $ws "E48" "  +2.95%  "
Set-TextCell $ws "E49" "  +0.51%  "
Set-TextCell $ws "D50" "0.999"
Set-TextCell $ws "E50" "  -0.21%  "
Set-TextCell $ws "D51" "134.30"
Set-TextCell $ws "E51" "  +1.12%  "
